$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'29.977.60"
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').Value = "'1.884.76"
$ws.Range('D3').Style = 'Normal'
$ws.Range('D4').Value = "'0.9995"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = "'0.7445"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.20%  '
$ws.Range('D6').Value = "'243.16"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.22%  '
$ws.Range('D7').Value = "'0.9999"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  +0.90%  '
$ws.Range('D9').Value = "'0.07236"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.37%  '
$ws.Range('E10').Value = '  -2.81%  '
$ws.Range('D11').Value = "'0.08346"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.20%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = "'2.006.31"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +4.62%  '
$ws.Range('B13').Value = 'Polygon'
$ws.Range('C13').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D13').Value = "'0.7571"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.98%  '
$ws.Range('D14').Value = "'5.421"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.89%  '
$ws.Range('D15').Value = "'92.76"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.05%  '
$ws.Range('E16').Value = '  +0.72%  '
$ws.Range('D17').Value = "'30.009.16"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.28%  '
$ws.Range('D18').Value = "'250.65"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.52%  '
$ws.Range('E19').Value = '  -1.15%  '
$ws.Range('D20').Value = "'0.000007868"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.64%  '
$ws.Range('D21').Value = "'2.220.27"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.42%  '
$ws.Range('D22').Value = "'1.000"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.09%  '
$ws.Range('D23').Value = "'8.055"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.34%  '
$ws.Range('D24').Value = "'0.9991"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.13%  '
$ws.Range('E25').Value = '  -3.95%  '
$ws.Range('D26').Value = "'9.321"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.80%  '
$ws.Range('D27').Value = "'165.69"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.60%  '
$ws.Range('D28').Value = "'18.75"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.20%  '
$ws.Range('E29').Value = '  +0.24%  '
$ws.Range('D30').Value = "'1.490"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.76%  '
$ws.Range('D31').Value = "'4.626"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.55%  '
$ws.Range('D32').Value = "'1.537"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.13%  '
$ws.Range('D33').Value = "'4.239"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.88%  '
$ws.Range('D34').Value = "'0.05377"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.43%  '
$ws.Range('E35').Value = '  +0.91%  '
$ws.Range('D36').Value = "'0.7592"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.68%  '
$ws.Range('D37').Value = "'0.9991"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.18%  '
$ws.Range('D38').Value = "'2.708"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.40%  '
$ws.Range('D39').Value = "'0.01971"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.04%  '
$ws.Range('D41').Value = "'0.4575"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.20%  '
$ws.Range('D42').Value = "'1.106.36"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.33%  '
$ws.Range('D43').Value = "'6.083"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.03%  '
$ws.Range('D44').Value = "'73.00"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.39%  '
$ws.Range('D45').Value = "'0.8657"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.63%  '
$ws.Range('D46').Value = "'104.64"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.55%  '
$ws.Range('E47').Value = '  +0.06%  '
$ws.Range('E48').Value = '  +0.04%  '
$ws.Range('D49').Value = "'7.626"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.72%  '
$ws.Range('D50').Value = "'2.076.39"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.25%  '
$ws.Range('D51').Value = "'9.543"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.24%  '
